$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "ignore" marker from E9 up to E2
$ws.Range("E9").ClearContents()
$ws.Range("E2").Value = "ignore"

# Update the active selection to E9 (reflecting the last selected cell)
$ws.Range("E9").Select()
